# Updated symbol list on Thu Dec 15 05:55:39 UTC 2022 with GitHub Actions
# Applies the price/coin-listing refresh described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Text)
    # Force Excel to store the value as TEXT (matches the workbook's
    # inline-string cells), even when the string looks like a number
    # (e.g. "264.72"). A leading apostrophe is Excel's own "treat as
    # text" marker; without it, numeric-looking strings silently become
    # numeric cells (and lose meaningful trailing zeros, e.g. "0.06150").
    # Resetting the Style back to "Normal" afterwards strips the
    # quote-prefix formatting Excel auto-applies, so the cell's visible
    # style stays identical to before the edit.
    $Range.Value = "'" + $Text
    $Range.Style = "Normal"
}

# --- Column D (Price) single-value updates -------------------------------
Set-TextValue $ws.Range("D2")  "264.72"
Set-TextValue $ws.Range("D3")  "22.72"
Set-TextValue $ws.Range("D4")  "6.289"
Set-TextValue $ws.Range("D5")  "0.06150"
Set-TextValue $ws.Range("D7")  "6.690"
Set-TextValue $ws.Range("D8")  "1.346"
Set-TextValue $ws.Range("D9")  "0.8304"
Set-TextValue $ws.Range("D10") "0.01355"
Set-TextValue $ws.Range("D11") "0.1579"
Set-TextValue $ws.Range("D12") "0.08205"
Set-TextValue $ws.Range("D14") "0.03134"

# --- Rows 15-26: coin list reshuffled (ProBitToken moved from rank 14 to
#     rank 25; every other row shifted up one slot) plus new prices -------
$rows = @(
    @{ Row = 15; B = "BitMartToken";            C = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx";      D = "0.09252";   E = "14BitMartTokenBMX" },
    @{ Row = 16; B = "MCDex";                   C = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb";                 D = "3.928";     E = "15MCDexMCB" },
    @{ Row = 17; B = "BitForexToken";           C = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf";      D = "0.001708";  E = "16BitForexTokenBF" },
    @{ Row = 18; B = "CoinExToken";             C = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet";       D = "0.04879";   E = "17CoinExTokenCET" },
    @{ Row = 19; B = "TigerCash";               C = "https://coinranking.com/coin/6hIn06L2+tigercash-tch";              D = "0.006229";  E = "18TigerCashTCH" },
    @{ Row = 20; B = "HotbitToken";             C = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb";        D = "0.005272";  E = "19HotbitTokenHTB" },
    @{ Row = 21; B = "BitKan";                  C = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan";            D = "0.001090";  E = "20BitKanKAN" },
    @{ Row = 22; B = "NitroEx";                 C = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx";             D = "0.0001500"; E = "21NitroExNTX" },
    @{ Row = 23; B = "LEO";                     C = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo";                D = "3.767";     E = "22LEOLEO" },
    @{ Row = 24; B = "BTSEToken";               C = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse";         D = "2.288";     E = "23BTSETokenBTSE" },
    @{ Row = 25; B = "BitpandaEcosystemToken";  C = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"; D = "0.3377";    E = "24BitpandaEcosystemTokenBEST" },
    @{ Row = 26; B = "ProBitToken";             C = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob";           D = "0.1227";    E = "25ProBitTokenPROB" }
)

foreach ($r in $rows) {
    $ws.Range("B" + $r.Row).Value = $r.B
    $ws.Range("C" + $r.Row).Value = $r.C
    Set-TextValue $ws.Range("D" + $r.Row) $r.D
    $ws.Range("E" + $r.Row).Value = $r.E
}

# --- Remaining column D single-value updates ------------------------------
Set-TextValue $ws.Range("D40") "0.04614"
Set-TextValue $ws.Range("D41") "0.006937"
Set-TextValue $ws.Range("D42") "0.1138"
Set-TextValue $ws.Range("D43") "0.003400"
Set-TextValue $ws.Range("D44") "0.01069"
Set-TextValue $ws.Range("D45") "0.00006163"
Set-TextValue $ws.Range("D48") "0.1958"
Set-TextValue $ws.Range("D49") "0.00002100"
